$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema5a"
$ws.Range("C2").Value = "Met"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2075446666666666
$ws.Range("H2").Value = 0.6226339999999999
$ws.Range("I2").Value = 0.01336584681749163
$ws.Range("J2").Value = 0.01405330443738086
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.912114666666667
$ws.Range("N2").Value = 8.736344000000001
$ws.Range("O2").Value = 0.04564398277650125
$ws.Range("P2").Value = 0.06163513710720567
$ws.Range("Q2").Value = 0.6043938677884444
$ws.Range("R2").Value = 5.439544810096
$ws.Range("S2").Value = 0.0006100704819309421
$ws.Range("T2").Value = 0.000866177345807271

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema5a"
$ws.Range("C3").Value = "Met"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2075446666666666
$ws.Range("H3").Value = 0.6226339999999999
$ws.Range("I3").Value = 0.01336584681749163
$ws.Range("J3").Value = 0.01405330443738086
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8623146666666667
$ws.Range("N3").Value = 2.586944
$ws.Range("O3").Value = 0.01351577128599483
$ws.Range("P3").Value = 0.01825095808139687
$ws.Range("Q3").Value = 0.1789688100551111
$ws.Range("R3").Value = 1.610719290496
$ws.Range("S3").Value = 0.0001806497286288588
$ws.Range("T3").Value = 0.0002564862701917466

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema5a"
$ws.Range("C4").Value = "Met"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2075446666666666
$ws.Range("H4").Value = 0.6226339999999999
$ws.Range("I4").Value = 0.01336584681749163
$ws.Range("J4").Value = 0.01405330443738086
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.668087
$ws.Range("N4").Value = 23.004261
$ws.Range("O4").Value = 0.12018827244785
$ws.Range("P4").Value = 0.1622956674765719
$ws.Range("Q4").Value = 1.591470560386
$ws.Range("R4").Value = 14.323235043474
$ws.Range("S4").Value = 0.001606418038796912
$ws.Range("T4").Value = 0.002280790423916196

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sema5a"
$ws.Range("C5").Value = "Met"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2075446666666666
$ws.Range("H5").Value = 0.6226339999999999
$ws.Range("I5").Value = 0.01336584681749163
$ws.Range("J5").Value = 0.01405330443738086
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 49.6589625
$ws.Range("N5").Value = 99.317925
$ws.Range("O5").Value = 0.7783460091712006
$ws.Range("P5").Value = 0.7006905777265834
$ws.Range("Q5").Value = 10.306452819075
$ws.Range("R5").Value = 61.83871691444999
$ws.Range("S5").Value = 0.0104032535295882
$ws.Range("T5").Value = 0.009847018005195952

$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Sema5a"
$ws.Range("C6").Value = "Met"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2075446666666666
$ws.Range("H6").Value = 0.6226339999999999
$ws.Range("I6").Value = 0.01336584681749163
$ws.Range("J6").Value = 0.01405330443738086
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.699147
$ws.Range("N6").Value = 8.097441
$ws.Range("O6").Value = 0.04230596431845346
$ws.Range("P6").Value = 0.05712765960824214
$ws.Range("Q6").Value = 0.5601935643993332
$ws.Range("R6").Value = 5.041742079593999
$ws.Range("S6").Value = 0.0005654550385467157
$ws.Range("T6").Value = 0.0008028323922696925

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema5a"
$ws.Range("C7").Value = "Met"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 13.025931
$ws.Range("H7").Value = 39.077793
$ws.Range("I7").Value = 0.8388680913725347
$ws.Range("J7").Value = 0.8820143483490313
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.912114666666667
$ws.Range("N7").Value = 8.736344000000001
$ws.Range("O7").Value = 0.04564398277650125
$ws.Range("P7").Value = 0.06163513710720567
$ws.Range("Q7").Value = 37.93300471208801
$ws.Range("R7").Value = 341.397042408792
$ws.Range("S7").Value = 0.03828928071436445
$ws.Range("T7").Value = 0.05436307529101521

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sema5a"
$ws.Range("C8").Value = "Met"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 13.025931
$ws.Range("H8").Value = 39.077793
$ws.Range("I8").Value = 0.8388680913725347
$ws.Range("J8").Value = 0.8820143483490313
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8623146666666667
$ws.Range("N8").Value = 2.586944
$ws.Range("O8").Value = 0.01351577128599483
$ws.Range("P8").Value = 0.01825095808139687
$ws.Range("Q8").Value = 11.232451348288
$ws.Range("R8").Value = 101.092062134592
$ws.Range("S8").Value = 0.01133794926211019
$ws.Range("T8").Value = 0.01609760689890875

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sema5a"
$ws.Range("C9").Value = "Met"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 13.025931
$ws.Range("H9").Value = 39.077793
$ws.Range("I9").Value = 0.8388680913725347
$ws.Range("J9").Value = 0.8820143483490313
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.668087
$ws.Range("N9").Value = 23.004261
$ws.Range("O9").Value = 0.12018827244785
$ws.Range("P9").Value = 0.1622956674765719
$ws.Range("Q9").Value = 99.883972163997
$ws.Range("R9").Value = 898.955749475973
$ws.Range("S9").Value = 0.1008221067136901
$ws.Range("T9").Value = 0.1431471073892197

$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Sema5a"
$ws.Range("C10").Value = "Met"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 13.025931
$ws.Range("H10").Value = 39.077793
$ws.Range("I10").Value = 0.8388680913725347
$ws.Range("J10").Value = 0.8820143483490313
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 49.6589625
$ws.Range("N10").Value = 99.317925
$ws.Range("O10").Value = 0.7783460091712006
$ws.Range("P10").Value = 0.7006905777265834
$ws.Range("Q10").Value = 646.8542190565876
$ws.Range("R10").Value = 3881.125314339525
$ws.Range("S10").Value = 0.6529296311408744
$ws.Range("T10").Value = 0.6180191433078187

$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Sema5a"
$ws.Range("C11").Value = "Met"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 13.025931
$ws.Range("H11").Value = 39.077793
$ws.Range("I11").Value = 0.8388680913725347
$ws.Range("J11").Value = 0.8820143483490313
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.699147
$ws.Range("N11").Value = 8.097441
$ws.Range("O11").Value = 0.04230596431845346
$ws.Range("P11").Value = 0.05712765960824214
$ws.Range("Q11").Value = 35.158902580857
$ws.Range("R11").Value = 316.430123227713
$ws.Range("S11").Value = 0.03548912354149561
$ws.Range("T11").Value = 0.05038741546206896

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Sema5a"
$ws.Range("C12").Value = "Met"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.2787875
$ws.Range("H12").Value = 4.557575
$ws.Range("I12").Value = 0.146753588727638
$ws.Range("J12").Value = 0.1028677986926446
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.912114666666667
$ws.Range("N12").Value = 8.736344000000001
$ws.Range("O12").Value = 0.04564398277650125
$ws.Range("P12").Value = 0.06163513710720567
$ws.Range("Q12").Value = 6.636090500966668
$ws.Range("R12").Value = 39.8165430058
$ws.Range("S12").Value = 0.006698418276274056
$ws.Range("T12").Value = 0.006340270876337582

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Sema5a"
$ws.Range("C13").Value = "Met"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.2787875
$ws.Range("H13").Value = 4.557575
$ws.Range("I13").Value = 0.146753588727638
$ws.Range("J13").Value = 0.1028677986926446
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.8623146666666667
$ws.Range("N13").Value = 2.586944
$ws.Range("O13").Value = 0.01351577128599483
$ws.Range("P13").Value = 0.01825095808139687
$ws.Range("Q13").Value = 1.965031883466667
$ws.Range("R13").Value = 11.7901913008
$ws.Range("S13").Value = 0.001983487940641704
$ws.Range("T13").Value = 0.001877435881865028

$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Sema5a"
$ws.Range("C14").Value = "Met"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.2787875
$ws.Range("H14").Value = 4.557575
$ws.Range("I14").Value = 0.146753588727638
$ws.Range("J14").Value = 0.1028677986926446
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 7.668087
$ws.Range("N14").Value = 23.004261
$ws.Range("O14").Value = 0.12018827244785
$ws.Range("P14").Value = 0.1622956674765719
$ws.Range("Q14").Value = 17.4739408045125
$ws.Range("R14").Value = 104.843644827075
$ws.Range("S14").Value = 0.01763806030469707
$ws.Range("T14").Value = 0.01669499805066839

$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Sema5a"
$ws.Range("C15").Value = "Met"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.2787875
$ws.Range("H15").Value = 4.557575
$ws.Range("I15").Value = 0.146753588727638
$ws.Range("J15").Value = 0.1028677986926446
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 49.6589625
$ws.Range("N15").Value = 99.317925
$ws.Range("O15").Value = 0.7783460091712006
$ws.Range("P15").Value = 0.7006905777265834
$ws.Range("Q15").Value = 113.1622230079687
$ws.Range("R15").Value = 452.648892031875
$ws.Range("S15").Value = 0.1142250701177087
$ws.Range("T15").Value = 0.07207849729541102

$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Sema5a"
$ws.Range("C16").Value = "Met"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.2787875
$ws.Range("H16").Value = 4.557575
$ws.Range("I16").Value = 0.146753588727638
$ws.Range("J16").Value = 0.1028677986926446
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.699147
$ws.Range("N16").Value = 8.097441
$ws.Range("O16").Value = 0.04230596431845346
$ws.Range("P16").Value = 0.05712765960824214
$ws.Range("Q16").Value = 6.1507824442625
$ws.Range("R16").Value = 36.904694665575
$ws.Range("S16").Value = 0.006208552088316445
$ws.Range("T16").Value = 0.005876596588362576

$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Sema5a"
$ws.Range("C17").Value = "Met"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.01572166666666667
$ws.Range("H17").Value = 0.047165
$ws.Range("I17").Value = 0.001012473082335678
$ws.Range("J17").Value = 0.001064548520943392
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 2.912114666666667
$ws.Range("N17").Value = 8.736344000000001
$ws.Range("O17").Value = 0.04564398277650125
$ws.Range("P17").Value = 0.06163513710720567
$ws.Range("Q17").Value = 0.04578329608444445
$ws.Range("R17").Value = 0.41204966476
$ws.Range("S17").Value = 0.000046213303931800841844665384
$ws.Range("T17").Value = 0.000065613594045619007505236453

$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("B18").Value = "Sema5a"
$ws.Range("C18").Value = "Met"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.01572166666666667
$ws.Range("H18").Value = 0.047165
$ws.Range("I18").Value = 0.001012473082335678
$ws.Range("J18").Value = 0.001064548520943392
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.8623146666666667
$ws.Range("N18").Value = 2.586944
$ws.Range("O18").Value = 0.01351577128599483
$ws.Range("P18").Value = 0.01825095808139687
$ws.Range("Q18").Value = 0.01355702375111111
$ws.Range("R18").Value = 0.12201321376
$ws.Range("S18").Value = 0.000013684354614075240784616151
$ws.Range("T18").Value = 0.000019429030431350889949549191

$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("B19").Value = "Sema5a"
$ws.Range("C19").Value = "Met"
$ws.Range("D19").Value = "Inflammatory-Mac"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.01572166666666667
$ws.Range("H19").Value = 0.047165
$ws.Range("I19").Value = 0.001012473082335678
$ws.Range("J19").Value = 0.001064548520943392
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 7.668087
$ws.Range("N19").Value = 23.004261
$ws.Range("O19").Value = 0.12018827244785
$ws.Range("P19").Value = 0.1622956674765719
$ws.Range("Q19").Value = 0.120555107785
$ws.Range("R19").Value = 1.084995970065
$ws.Range("S19").Value = 0.0001216873906658749
$ws.Range("T19").Value = 0.0001727716127677053

$ws.Range("A20").Value = "Resolving-Mac"
$ws.Range("B20").Value = "Sema5a"
$ws.Range("C20").Value = "Met"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.01572166666666667
$ws.Range("H20").Value = 0.047165
$ws.Range("I20").Value = 0.001012473082335678
$ws.Range("J20").Value = 0.001064548520943392
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 49.6589625
$ws.Range("N20").Value = 99.317925
$ws.Range("O20").Value = 0.7783460091712006
$ws.Range("P20").Value = 0.7006905777265834
$ws.Range("Q20").Value = 0.7807216554374999
$ws.Range("R20").Value = 4.684329932625
$ws.Range("S20").Value = 0.0007880543830292397
$ws.Range("T20").Value = 0.0007459191181578054

$ws.Range("A21").Value = "Resolving-Mac"
$ws.Range("B21").Value = "Sema5a"
$ws.Range("C21").Value = "Met"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.01572166666666667
$ws.Range("H21").Value = 0.047165
$ws.Range("I21").Value = 0.001012473082335678
$ws.Range("J21").Value = 0.001064548520943392
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 2.699147
$ws.Range("N21").Value = 8.097441
$ws.Range("O21").Value = 0.04230596431845346
$ws.Range("P21").Value = 0.05712765960824214
$ws.Range("Q21").Value = 0.04243508941833333
$ws.Range("R21").Value = 0.381915804765
$ws.Range("S21").Value = 0.000042833650094687806826176252
$ws.Range("T21").Value = 0.000060815165540911747977127072
